$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily level")

$ws.Range("E1").Value = "Fruit Servings(Cal)"
$ws.Range("F1").Value = "Vegetable Servings(Cal)"

$ws.Range("G1").Select()
